# Generate Report for Handback
# Row 8 (the c0484e70-... file) on the zh-cn and de-de sheets currently
# reuses row 7's "Correspond Handoff Datetime" / "Correspond Handback
# DateTime" strings. This fills in the real, distinct timestamps for
# row 8's handoff/handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D8").Value = "2016-03-08 02:17:17"
$wsZhCn.Range("G8").Value = "2016-03-08 02:17:59"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D8").Value = "2016-03-08 02:17:25"
$wsDeDe.Range("G8").Value = "2016-03-08 02:18:12"
